# Generate Report for handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) values for the row that
# corresponds to the second handback file (row 3) on both the "zh-cn"
# and "de-de" report sheets, reflecting the newly generated handback
# report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-08 07:39:12"
$wsZhCn.Range("G3").Value = "2016-01-08 07:39:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-08 07:39:22"
$wsDeDe.Range("G3").Value = "2016-01-08 07:40:13"
